$d = $word.ActiveDocument

# Replace paragraph text directly (Range.Text) rather than using
# Find/Replace, since Find.Execute's replacement text passes through
# AutoFormat/AutoCorrect and silently turns straight apostrophes into
# curly ones. Assigning Range.Text keeps the existing trailing paragraph
# mark intact, so no trailing `r should be appended here.
$d.Paragraphs(2).Range.Text = "**PRD Draft: Vibe Coders Prompt Management & Enhancement Tool**"
$d.Paragraphs(3).Range.Text = "**1. Problem Summary:** Vibe Coders, users of AI coding tools, experience significant challenges related to prompt management, quality, and consistency. These include wasted time due to repetitive prompt creation, difficulty managing prompt libraries across platforms, inconsistent LLM outputs, challenges in controlling LLM tone and behavior, and concerns about AI hallucinations and data accuracy.  Furthermore, a lack of community support, especially for Spanish-speaking users, hinders effective collaboration and knowledge sharing around prompt engineering best practices."
$d.Paragraphs(4).Range.Text = "**2. Why This Problem Matters:**  These problems directly impact Vibe Coders' productivity and the quality of their code.  Wasted time on repetitive tasks reduces coding efficiency and slows down project completion. Inconsistent LLM outputs lead to debugging delays and increased development costs.  The lack of a centralized and collaborative prompt management system hinders knowledge sharing and slows down team development cycles.  For businesses using Vibe Coders' work, these inefficiencies translate into missed deadlines, higher development costs, and potentially compromised software quality.  The absence of adequate support in Spanish further limits accessibility and inclusivity within the community."
$d.Paragraphs(5).Range.Text = "**3. Potential Solution Overview:** We propose developing a comprehensive prompt management and enhancement tool integrated directly into existing AI coding tools or as a standalone application. This tool will address prompt organization, quality control, community collaboration, and mitigation of LLM limitations.  The solution will include features to streamline prompt creation, storage, version control, and sharing, along with integrated tools for prompt engineering guidance and hallucination detection."
$d.Paragraphs(7).Range.Text = "* **Centralized Prompt Library:**  Users can create, save, organize, and version-control their prompts in a single, searchable library accessible across platforms. This solves the problem of disorganized prompt libraries and repetitive typing."
$d.Paragraphs(8).Range.Text = "* **Prompt Quality Assurance Tools:**  The tool will integrate basic prompt engineering guidance, character count monitoring, and a simple hallucination detection mechanism (e.g., flagging potentially inaccurate outputs). This addresses issues with prompt quality and LLM reliability."
$d.Paragraphs(9).Range.Text = "* **Community Collaboration Features (Beta - Spanish Support):**  A basic forum/discussion board will be implemented, initially focusing on Spanish-speaking users, to facilitate knowledge sharing and support around prompt engineering and best practices.  This addresses the identified need for community support and inclusivity."
$d.Paragraphs(10).Range.Text = "* **Prompt Templates & Snippets:** Pre-built prompt templates and reusable code snippets will accelerate prompt creation for common coding tasks. This improves efficiency and reduces repetitive work."
$d.Paragraphs(11).Range.Text = "* **Tone & Style Control (Basic):**  Users will have basic controls to influence the tone and style of LLM responses (e.g., formal vs. informal). This addresses challenges in controlling LLM behavior."
$d.Paragraphs(13).Range.Text = "* **User Interviews (1 week):** Conduct 5-7 user interviews with Vibe Coders to validate the MVP features and gather further insights into their workflow and needs.  Focus on Spanish-speaking users to understand their specific requirements."
$d.Paragraphs(14).Range.Text = "* **Prototype Development (2 weeks):** Develop a low-fidelity prototype of the centralized prompt library and basic prompt quality assurance features."
$d.Paragraphs(15).Range.Text = "* **Sprint Planning (1 week):** Based on user feedback and prototype testing, finalize the MVP feature set and create a detailed sprint plan for development."
$d.Paragraphs(16).Range.Text = "* **A/B Testing (Ongoing):** After launch, conduct A/B testing to compare the effectiveness of the new tool against existing workflows."

# Append a new closing paragraph at the end of the document
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter("`r`nThis PRD provides a starting point for development.  Further refinement will be necessary based on user research and ongoing feedback.")
